$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply quotePrefix style to B66:B68 (re-enter "Vrij" with a leading apostrophe) ---
$ws.Range("B66").Value = "'Vrij"
$ws.Range("B67").Value = "'Vrij"
$ws.Range("B68").Value = "'Vrij"

# --- Week 13 (rows 66-71) hours log entries ---
# Numeric / shared-cell values first
$ws.Range("B69").Value = 2
$ws.Range("C69").Value = 2
$ws.Range("G69").Value = 2
$ws.Range("I69").Value = 2
$ws.Range("K69").Value = 2
$ws.Range("M69").Value = 2
$ws.Range("O69").Value = 2

$ws.Range("B70").Value = 4
$ws.Range("C70").Value = 4
$ws.Range("G70").Value = 4
$ws.Range("I70").Value = 4
$ws.Range("K70").Value = 4
$ws.Range("M70").Value = 3.15
$ws.Range("O70").Value = 3.75
$ws.Range("P70").Value = "kwartier te laat"

# New remark strings, entered in the same order they were first introduced
# (matches the shared-string table ordering: N70, Q69, Q70, F69, F70)
$ws.Range("N70").Value = "drie kwartier te laat"
$ws.Range("Q69").Value = "Harold ziek"
$ws.Range("Q70").Value = "Harold ziek, danial kwartier te laat, paco drie kwartier te laat"
$ws.Range("F69").Value = "ziek"
$ws.Range("F70").Value = "ziek, zou proberen thuis te werken"

# Row 71 totals: B71 gains a real SUM formula like its siblings
$ws.Range("B71").Formula = "=SUM(B69:B70)"

# --- Outline / grouping changes ---
# Collapse (hide) the week 8 and week 9 detail rows
$ws.Range("A41:A47").EntireRow.Hidden = $true
$ws.Range("A49:A55").EntireRow.Hidden = $true
# Expand (show) the week 13 detail rows
$ws.Range("A73:A79").EntireRow.Hidden = $false

# Un-hide column P (16), which stays part of the outline group but is visible now
$ws.Columns.Item(16).Hidden = $false

# --- Sheet view: clear the frozen/scrolled topLeftCell and move the selection to Q67 ---
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("Q67").Select()
